$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Emp Info")

# Existing students (rows 2-4) get re-inscribed: niveau goes 1 -> 2, type goes
# INSCRIPTION -> REINSCRIPTION.
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(2, 6).Value = "REINSCRIPTION"
$ws.Cells.Item(3, 6).Value = "REINSCRIPTION"
$ws.Cells.Item(4, 6).Value = "REINSCRIPTION"

# Row 3's student also gets a name tweak.
$ws.Cells.Item(3, 3).Value = "Benabbou2"
$ws.Cells.Item(3, 4).Value = "Oussama2"

# New students registered (rows 5-7), freshly inscribed.
$ws.Cells.Item(5, 1).Value = 103
$ws.Cells.Item(6, 1).Value = 104
$ws.Cells.Item(7, 1).Value = 105

$ws.Cells.Item(5, 2).Value = "A133341333"
$ws.Cells.Item(6, 2).Value = "B133341333"
$ws.Cells.Item(7, 2).Value = "C133341333"

$ws.Cells.Item(5, 3).Value = "Yamani"
$ws.Cells.Item(5, 4).Value = "Jamal"
$ws.Cells.Item(6, 3).Value = "Ferdous"
$ws.Cells.Item(6, 4).Value = "Kamal"
$ws.Cells.Item(7, 3).Value = "Touhami"
$ws.Cells.Item(7, 4).Value = "Badr"

$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(7, 5).Value = 1

$ws.Cells.Item(5, 6).Value = "INSCRIPTION"
$ws.Cells.Item(6, 6).Value = "INSCRIPTION"
$ws.Cells.Item(7, 6).Value = "INSCRIPTION"

# Match the final selection left behind in the sheet.
$ws.Range("F10").Select()
